{"js": "// Fill the first empty row of the time-tracking table with a new entry:\n// date 17.02.23, 1,5 hours, and a three-line description.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\n// Find the first row (after the header) whose cells are all empty text \u2014\n// that is the blank row the new entry belongs in.\nconst rowCount = table.rows.items.length;\nlet targetRowIndex = -1;\nfor (let i = 1; i < rowCount; i++) {\n  const dateCell = table.getCell(i, 0);\n  dateCell.load(\"body/text\");\n  const hoursCell = table.getCell(i, 1);\n  hoursCell.load(\"body/text\");\n  const descCell = table.getCell(i, 2);\n  descCell.load(\"body/text\");\n  await context.sync();\n\n  if (\n    dateCell.body.text.trim() === \"\" &&\n    hoursCell.body.text.trim() === \"\" &&\n    descCell.body.text.trim() === \"\"\n  ) {\n    targetRowIndex = i;\n    break;\n  }\n}\n\nif (targetRowIndex === -1) {\n  throw new Error(\"No empty row found to fill in.\");\n}\n\nconst dateCell = table.getCell(targetRowIndex, 0);\ndateCell.body.insertText(\"17.02.23\", Word.InsertLocation.replace);\n\nconst hoursCell = table.getCell(targetRowIndex, 1);\nhoursCell.body.insertText(\"1,5\", Word.InsertLocation.replace);\n\nconst descCell = table.getCell(targetRowIndex, 2);\ndescCell.body.insertText(\n  \"Kirjautumissivu k\u00e4ytt\u00e4m\u00e4\u00e4n tietokantaa, \\n\" +\n    \"K\u00e4ytt\u00e4j\u00e4ntunnuksen muutos pieniksi kirjaimiksi,\\n\" +\n    \"K\u00e4ytt\u00e4j\u00e4sivun k\u00e4ytt\u00e4j\u00e4tunnus textbox read only\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Fill the first empty row of the time-tracking table with a new entry:\n# date 17.02.23, 1,5 hours, and a three-line description.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$targetRow = -1\nfor ($i = 2; $i -le $rowCount; $i++) {\n    $c1 = $t.Cell($i, 1).Range.Text -replace \"[\\r\\a]\", \"\"\n    $c2 = $t.Cell($i, 2).Range.Text -replace \"[\\r\\a]\", \"\"\n    $c3 = $t.Cell($i, 3).Range.Text -replace \"[\\r\\a]\", \"\"\n    if ($c1 -eq \"\" -and $c2 -eq \"\" -and $c3 -eq \"\") {\n        $targetRow = $i\n        break\n    }\n}\n\nif ($targetRow -eq -1) {\n    throw \"No empty row found to fill in.\"\n}\n\n$dateCell = $t.Cell($targetRow, 1)\n$dateCell.Range.Text = \"17.02.23\"\n\n$hoursCell = $t.Cell($targetRow, 2)\n$hoursCell.Range.Text = \"1,5\"\n\n$descCell = $t.Cell($targetRow, 3)\n$descCell.Range.Text = \"Kirjautumissivu k\u00e4ytt\u00e4m\u00e4\u00e4n tietokantaa, `rK\u00e4ytt\u00e4j\u00e4ntunnuksen muutos pieniksi kirjaimiksi,`rK\u00e4ytt\u00e4j\u00e4sivun k\u00e4ytt\u00e4j\u00e4tunnus textbox read only\"\n"}
